$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "AvailableColumns" list in column C (rows 2-51 originally) is an
# alphabetically-sorted list of column names. This commit removes the entry
# "Request Addl %" and introduces three new entries: "Global List Price(USD)",
# "Global Net Price(USD)" and "Request Total %". Re-sorting alphabetically
# (case-insensitive) after that change produces the 52-item list below
# (rows 2-14 are unaffected by the change, so only rows 15-53 need updating,
# plus two brand new rows for the list growing by two net entries).

$ws.Cells.Item(15,3).Value2 = "Global List Price(USD)"
$ws.Cells.Item(16,3).Value2 = "Global Net Price(USD)"
$ws.Cells.Item(17,3).Value2 = "Group Name"
$ws.Cells.Item(18,3).Value2 = "Group Number"
$ws.Cells.Item(19,3).Value2 = "HP Cost"
$ws.Cells.Item(20,3).Value2 = "Item Level"
$ws.Cells.Item(21,3).Value2 = "Item Number"
$ws.Cells.Item(22,3).Value2 = "List Price"
$ws.Cells.Item(23,3).Value2 = "Manufacturing Product ID"
$ws.Cells.Item(24,3).Value2 = "Master Contract ID"
$ws.Cells.Item(25,3).Value2 = "MCC"
$ws.Cells.Item(26,3).Value2 = "My Empowerment Disc %"
$ws.Cells.Item(27,3).Value2 = "Open Market"
$ws.Cells.Item(28,3).Value2 = "Option"
$ws.Cells.Item(29,3).Value2 = "PA Expiration Date"
$ws.Cells.Item(30,3).Value2 = "PA Number"
$ws.Cells.Item(31,3).Value2 = "Preferred Disc %"
$ws.Cells.Item(32,3).Value2 = "Preferred Suppliers"
$ws.Cells.Item(33,3).Value2 = "Price Structure"
$ws.Cells.Item(34,3).Value2 = "Pricing Source"
$ws.Cells.Item(35,3).Value2 = "Product Class"
$ws.Cells.Item(36,3).Value2 = "Product Description"
$ws.Cells.Item(37,3).Value2 = "Product Line"
$ws.Cells.Item(38,3).Value2 = "Quantity On Hand"
$ws.Cells.Item(39,3).Value2 = "Recycling Fee"
$ws.Cells.Item(40,3).Value2 = "Reference Price(US$)"
$ws.Cells.Item(41,3).Value2 = "Request Net"
$ws.Cells.Item(42,3).Value2 = "Request Total %"
$ws.Cells.Item(43,3).Value2 = "Serial Code"
$ws.Cells.Item(44,3).Value2 = "Solution ID"
$ws.Cells.Item(45,3).Value2 = "Source ID"
$ws.Cells.Item(46,3).Value2 = "Supplier Code"
$ws.Cells.Item(47,3).Value2 = "System ID"
$ws.Cells.Item(48,3).Value2 = "System Name"
$ws.Cells.Item(49,3).Value2 = "TAA Compliance"
$ws.Cells.Item(50,3).Value2 = "Total Requested Discount"
$ws.Cells.Item(51,3).Value2 = "Unit Price"
$ws.Cells.Item(52,3).Value2 = "Unit Weight"
$ws.Cells.Item(53,3).Value2 = "Warranty Code"

# Column C is a bit wider now to fit the new, longer entries.
$ws.Columns.Item(3).ColumnWidth = 26.3

# Update the view: selection moved to F11 and the previous scrolled
# position (topLeftCell) is cleared.
$ws.Range("F11").Select()
